$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix wording in the existing "Lösungsvorschlag" cell (B5): add "das" before "Telemetrie"
$ws.Range("B5").Value = "Windows 10 Enterprise benutzen und somit auch das Telemetrie Problem lösen"

# Add a new row 6 with a "Begründung" (justification) entry
$ws.Range("A6").Value = "Begründung"
$ws.Range("B6").Value = "Mit geeignter Version Problemen aus dem Weg gehen"

# Widen column B slightly to accommodate the new content (closest achievable width)
$ws.Range("B1").ColumnWidth = 70.666666666667

# Move the active selection, as it ended up after the last edits
[void]$ws.Range("B10").Select()
